# Actualización 10 de Mayo
$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 1P" ---
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D2").Value = 3
$ws1.Range("F2").Value = 22
$ws1.Range("G2").Value = 88
$ws1.Range("H2").Value = 7.8

$ws1.Range("D4").Value = 5
$ws1.Range("F4").Value = 16
$ws1.Range("G4").Value = 76.19
$ws1.Range("H4").Value = 7.8

# --- Sheet "Estadisticos 2P" ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("D2").Value = 4
$ws2.Range("E2").Value = 1
$ws2.Range("F2").Value = 21
$ws2.Range("G2").Value = 84
$ws2.Range("H2").Value = 7.9

$ws2.Range("D3").Value = 1
$ws2.Range("E3").Value = 0
$ws2.Range("F3").Value = 24
$ws2.Range("G3").Value = 96
$ws2.Range("H3").Value = 7.8

$ws2.Range("D4").Value = 5
$ws2.Range("E4").Value = 0
$ws2.Range("F4").Value = 16
$ws2.Range("G4").Value = 76.19
$ws2.Range("H4").Value = 7.8

# --- Sheet "Estadisticos Final" ---
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D2").Value = 3
$ws3.Range("F2").Value = 22
$ws3.Range("G2").Value = 88
$ws3.Range("H2").Value = 7.5

$ws3.Range("H3").Value = 7.6

$ws3.Range("D4").Value = 5
$ws3.Range("F4").Value = 16
$ws3.Range("G4").Value = 76.19
$ws3.Range("H4").Value = 7.7
